# Append new scrape run (2026-01-29 06:43 JST) to the "ランサーズ" sheet.
# New unique postings are prepended; two postings that are still active
# (5481153, 5481091) are kept/updated in place; everything else falls off
# the bottom of the (now shorter) list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("ランサーズ")

# --- 1. Clear all existing hyperlinks (and their relationships) up front so
#        we can rebuild a clean, correctly-ordered set at the end. Row
#        deletion below does not keep the hyperlinks collection in sync, so
#        doing this first avoids stale / duplicated entries. ---
$ws.Hyperlinks.Delete()

# --- 2. Drop the rows that fall off the bottom of the list (old rows 7-11).
#        This also shifts the dimension down to A1:H6. ---
$ws.Range("A7:A11").EntireRow.Delete()

# --- 3. Overwrite rows 2-6 with the latest scrape snapshot. ---

$ws.Range("A2").Value = "2026-01-29 06:43:37"
$ws.Range("B2").Value = "【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5473940"
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = "◆開発 ◇業務改善"

$ws.Range("A3").Value = "2026-01-29 06:43:37"
$ws.Range("B3").Value = "【フルリモート】セールスフォース開発案件の設計〜実装"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5481435"
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = "◆開発"

$ws.Range("A4").Value = "2026-01-29 06:43:37"
$ws.Range("B4").Value = "【長期・フルリモート】Webサイト修正業務(Git使用/Cursor環境)|安定稼働できる方歓迎"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5481153"
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = "◇サイト"

$ws.Range("A5").Value = "2026-01-29 06:43:37"
$ws.Range("B5").Value = "【ストレスチェックシステム】運用コスト削減と個人情報保護"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5481296"
$ws.Range("G5").Value = 33
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = "2026-01-29 06:43:37"
$ws.Range("B6").Value = "【短期/読み取り専用】既存システムに触れない行動ログ基盤の構築(Fintech系)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5481091"
$ws.Range("G6").Value = 33
$ws.Range("H6").ClearContents()

# --- 4. Re-create the hyperlinks on the URL column for the 5 remaining rows. ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473940")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5481435")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5481153")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5481296")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5481091")

# Adding hyperlinks one-by-one stamps a brand-new "Hyperlink" cell style onto
# each cell instead of reusing the workbook's existing style slot. Re-apply
# the named style so the cells collapse back onto the original style index.
$ws.Range("F2:F6").Style = "Normal"
$ws.Range("F2:F6").Style = "Hyperlink"

# --- 5. Column width tweaks (B, D, H). ColumnWidth round-trips through
#        Excel's internal "characters" unit and comes back ~0.8333 wider
#        than what is stored in the saved file, so compensate here. ---
$offset = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 50 - $offset
$ws.Columns.Item(4).ColumnWidth = 28 - $offset
$ws.Columns.Item(8).ColumnWidth = 12 - $offset

Write-Host "done"
